$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B7").Value = 'The marking scheme is highly accurate and ready for use. The terminology matches current Hong Kong educational policy (VPET) and VTC’s institutional structure. No factual errors were found. This evaluation confirms that the marking scheme is factually accurate and aligns with the official standards and terminology of the Vocational Training Council (VTC) in Hong Kong. [vtc.edu.hk](https://vertexaisearch.cloud.google.com/grounding-api-redirect/AUZIYQGfYjWT3BilK0Z-osBdBJkrXPQjTZDsS96PhghiYYnUHdoFpqFUpVSQHZDE4dIKCXYhvEE8x1Z4-BwC-kdHBBBLMTaY8h9bLVoJhUOxj-2i2_Jt848kCldR_bb5ZQ3CyZuT4EUrO9Yj)'
$wsSummary.Range("B8").Value = '2026-01-07 14:20:19'
$wsVerif = $wb.Worksheets.Item("Content Verification")
$wsVerif.Range("C2").Value = 'Terminology is 100% correct. VTC is officially recognized as the largest VPET provider in HK. According to the Hong Kong Education Bureau (EDB) and VTC corporate documents, VPET stands for Vocational and Professional Education and Training. VTC is indeed the largest provider, serving over 200,000 students annually. [vtc.edu.hk](https://vertexaisearch.cloud.google.com/grounding-api-redirect/AUZIYQGXnV_z3VfqDjGu28vBBuPopdDzvAy0H5JnyvN-0fZvBU_rcHRuERR8w7TzXA5U3LzaFmmxiUo_5iKb8ch-sNZok_mgIJlO9QFl0r9JixDHfRw5-TS_gF_ma1Y=)'
$wsVerif.Range("D2").Value = 'The marking scheme is excellent. To make it even more robust, you could mention that VPET is recognized under the Hong Kong Qualifications Framework (HKQF) to emphasize its formal status.'
$wsVerif.Range("C3").Value = 'Correctly distinguishes between the sub-degree (IVE) and degree (THEi) levels. IVE (Hong Kong Institute of Vocational Education) is the primary member institution for Higher Diplomas. THEi (Technological and Higher Education Institute of Hong Kong) was specifically established to offer vocationally-oriented Bachelor’s Degrees. [apqn.org](https://vertexaisearch.cloud.google.com/grounding-api-redirect/AUZIYQHfJJY7QCN6DCC9qBFoBZZW_kLBax2iyr-W9BtD79JuasnZbUtE6JbYEUPxCxFoU0mVdszotrSNl4En9a-dNzkWcOQhSfLz2pgFUX6CfXoerxLfsR2t_yfQAx5v3xYLbAyqVgd2Ss3xHUNuOSemW_Gl73zrWrsHSqo=)'
$wsVerif.Range("D3").Value = 'The answer is correct. You might add that THEi degrees are often "Applied Degrees," which is a specific term currently promoted by the EDB.'
$wsVerif.Range("C4").Value = 'This is the official VTC motto. The explanation of theory vs. practice is accurate. VTC’s official branding frequently uses the phrase "Think and Do" (sometimes paired with "Hands-on, Minds-on"). It refers to the integration of theoretical knowledge with practical application. [vtc.edu.hk](https://vertexaisearch.cloud.google.com/grounding-api-redirect/AUZIYQHDbintfm2E7gNTfbWuVqyf9FvW7zrwNcMEOKd3Uo9xfb_TXpbT5trnUT3u5ZWbshmBj0eCp5Ig0B-4BUSk0JcyCpjzS3nP6Yj-C1bK2h0HYJJmKhoGD0v4WxTxT_2NaMn4OT3HbrpnRdl-5KlKPLv_ElwCA==)'
$wsVerif.Range("D4").Value = 'The marking breakdown is well-balanced. Ensure students understand that "Think" isn''t just "remembering facts" but "problem-solving and analytical thinking."'
$wsVerif.Range("C5").Value = 'DFS and DVE are the standard "bridge" programmes for S6 students. For S6 students who do not meet the "3322+2" (Degree) or "22222" (HD) DSE requirements, the Diploma of Foundation Studies (DFS) is the standard one-year progression route. The Diploma of Vocational Education (DVE) is also a valid answer, typically offered by Youth College. [vtc.edu.hk](https://vertexaisearch.cloud.google.com/grounding-api-redirect/AUZIYQFk0SBJGwIzx2bq4cHmucdpXKcFbf8_K2l1Uicysm9GGnfQmybL0cEF8kGz5lGHU8x1DeP1g0ZCMPsR_ARpYI8irjnn4Ycl11KBmNbimHQu3omhmH6Ekq_klGPT17wI714YDkfP_JODmtafU0GWRlPtZO4Jfk2DNRWJiiUltvp-pV7HIg==)'
$wsVerif.Range("D5").Value = 'The answer is accurate. Note that successful completion of DFS is considered equivalent to Level 2 in five DSE subjects (including English and Chinese) for civil service applications, which is a major "importance" factor you could include.'
$wsVerif.Range("C6").Value = 'Benefits listed (internships, equipment) are core pillars of VTC’s "Work-Integrated Learning." VTC operates under a "tripartite" model involving the government, industry, and the council. Work-Integrated Learning (WIL) and the VTC Earn & Learn Scheme are real-world examples of this collaboration. [vtc.edu.hk](https://vertexaisearch.cloud.google.com/grounding-api-redirect/AUZIYQHwIxQUfSgNHoK6cwuQKfqGy010xqQsQ0Xpits1ZLg3iuvuTwFIxoe0JCSkhIaDZYw3kuDzP5t7rLPalvtMDy5FtsLTy8JT3IHuhjvK707d5tKiAtjeJArzHq2QqSU58pHznLcPz0g3mEdH2pHlwQbJ7IJM0AoN7RyyCKR0rt3K0A1kzueNR4A6C_lJakU=)'
$wsVerif.Range("D6").Value = 'The examples provided (internships, equipment) are perfect. You could also suggest "Industry-endorsed curriculum" or "Professional recognition/certification" as alternative valid examples.'
